$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "call getStartParams()" test case (previously row 11) is being fixed and
# moved up to row 8. Rows that used to be 8, 9 and 10 shift down to 9, 10, 11.
# Capture the old content of rows 8, 9, 10 and 11 BEFORE overwriting anything.

$oldE8 = $ws.Range("E8").Value2
$oldG8 = $ws.Range("G8").Value2
$oldH8 = $ws.Range("H8").Value2

$oldE9 = $ws.Range("E9").Value2
$oldG9 = $ws.Range("G9").Value2
$oldH9 = $ws.Range("H9").Value2

$oldE10 = $ws.Range("E10").Value2
$oldG10 = $ws.Range("G10").Value2
$oldH10 = $ws.Range("H10").Value2

$oldE11 = $ws.Range("E11").Value2

# New combined steps / expected-results text for the fixed getStartParams test.
$newG8 = @'
wait(3);
validate1;
link_Click(system_test_link);
validate2;
SelectTestToRun(VT200_0963_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(3);
SelectTestToRun(VT200_0959_string);
ClickRunTest(runtest_top_xpath);
validate4;
ClickRunTest(runtest_bottom_xpath);
TaponGetStartparams;
validate5;
CheckUITextContains(?ParamsAreSet);
ClickUIButtonText(OK);
press_Key(Back);
'@

$newH8 = @'
validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=System JS Test
};
validate3
{
validate_Text_Exists=VT200-0963
};
validate4
{
validate_Text_Exists=VT200-0959
};
validate5
{
validate_App_Launched_Device=com.rhomobile.testapp
};
'@

# Row 8 becomes the getStartParams test case, with the new steps/results text.
$ws.Range("E8").Value = $oldE11
$ws.Range("G8").Value = $newG8
$ws.Range("H8").Value = $newH8

# Old row 8 content moves down to row 9.
$ws.Range("E9").Value = $oldE8
$ws.Range("G9").Value = $oldG8
$ws.Range("H9").Value = $oldH8

# Old row 9 content moves down to row 10.
$ws.Range("E10").Value = $oldE9
$ws.Range("G10").Value = $oldG9
$ws.Range("H10").Value = $oldH9

# Old row 10 content moves down to row 11.
$ws.Range("E11").Value = $oldE10
$ws.Range("G11").Value = $oldG10
$ws.Range("H11").Value = $oldH10

# Row heights follow the (longer) moved/updated content.
$ws.Rows.Item(8).RowHeight = 268.5
$ws.Rows.Item(9).RowHeight = 255.75
$ws.Rows.Item(10).RowHeight = 217.5

# Update the view state to match where the edit was made.
$ws.Range("H8").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
